# Fix the logo lettering position on slide 1 ("page not rendeinrg logo").
#
# The five single-letter text boxes that spell the logo lettering
# ("S","T","A","T","S") live inside the top-level group shape on slide 1
# and were nudged to a new position (shifted by +107325 EMU in x and
# +64395 EMU in y). The sixth label ("SMUSL") is left untouched.
#
# PowerPoint's Shape.Left / Shape.Top are expressed in points (1 pt =
# 12700 EMU) and are stored internally as 32-bit floats, so the literal
# point values below were chosen (via float32 round-trip search) so
# that re-serializing them back to EMU reproduces the exact target
# offsets from the authoritative edit.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$logoGroup = $s.Shapes.Item(1)

# name -> [Left(pt), Top(pt)]  (target EMU noted alongside for clarity)
$moves = @{
    "TextBox{S2}" = @(508.00079345703125, 59.75165557861328)  # -> x=6451610 y=758846
    "TextBox{T2}" = @(453.80938720703125, 59.75165557861328)  # -> x=5763379 y=758846
    "TextBox{A}"  = @(372.175537109375,   59.75165557861328)  # -> x=4726629 y=758846
    "TextBox{T1}" = @(323.7422180175781,  59.75165557861328)  # -> x=4111526 y=758846
    "TextBox{S1}" = @(273.25079345703125, 59.75165557861328)  # -> x=3470285 y=758846
}

foreach ($name in $moves.Keys) {
    $shape = $logoGroup.GroupItems.Item($name)
    $xy = $moves[$name]
    $shape.Left = $xy[0]
    $shape.Top = $xy[1]
}
